# Applies the "20 minutes pres" update: adds three new weekly weight-measurement
# columns (H, I, J) for dates 45378, 45385, 45393, mirroring the formatting of
# the existing column G.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of column G (header style + data style) into H:J so the
# new columns look identical to the existing weekly measurement columns.
$srcRange = $ws.Range("G1:G27")
$dstRange = $ws.Range("H1:J27")
$srcRange.Copy() | Out-Null
$dstRange.PasteSpecial(-4122) | Out-Null   # xlPasteFormats = -4122
$excel.CutCopyMode = 0

# New measurement dates (row 1 header, serial date numbers matching the
# weekly cadence already used in columns E/F/G).
$ws.Range("H1").Value = 45378
$ws.Range("I1").Value = 45385
$ws.Range("J1").Value = 45393

# New weight measurements for each animal row (2-27).
$data = @(
    @(2, 20.2, 20.6, 20.3),
    @(3, 22.3, 21.6, 22.3),
    @(4, 23.6, 24, 24.1),
    @(5, 23.8, 23.7, 23.2),
    @(6, 23.1, 23.9, 23.9),
    @(7, 20.4, 20.4, 20.9),
    @(8, 20.5, 21.2, 21.2),
    @(9, 21.8, 21.3, 21.5),
    @(10, 21.1, 21.6, 22),
    @(11, 21.6, 21.7, 23.6),
    @(12, 22.1, 24.4, 23.8),
    @(13, 22.2, 24.7, 22.7),
    @(14, 22.5, 23.9, 22.5),
    @(15, 20.9, 21, 20.5),
    @(16, 24.8, 24.2, 24.8),
    @(17, 23.2, 22.5, 23.4),
    @(18, 22, 21.6, 21.9),
    @(19, 21.8, 22, 22.3),
    @(20, 20, 20.2, 20.6),
    @(21, 21.8, 21.8, 22.8),
    @(22, 18.9, 19.7, 20.1),
    @(23, 20.2, 21.2, 20.7),
    @(24, 24.4, 24.5, 23.5),
    @(25, 24.3, 23, 25.2),
    @(26, 20.6, 20.7, 21.7),
    @(27, 22.1, 22.8, 23.7)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 8).Value = $row[1]
    $ws.Cells.Item($r, 9).Value = $row[2]
    $ws.Cells.Item($r, 10).Value = $row[3]
}

# Mirror the final selection state recorded by Excel after the paste
# (active cell E1, selected range E1:J27).
$ws.Range("E1:J27").Select() | Out-Null
